$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 62 col A: drop the trailing ellipsis
$ws.Range("A62").Value = "Results for this location"

# Row 65 (previously blank) gets new content
$ws.Range("A65").Value = "What to Expect at This Location"
$ws.Range("B65").Value = "Qué esperar en esta ubicación"

# New row 66
$ws.Range("A66").Value = "Getting results for your location…"
$ws.Range("B66").Value = "Obtener resultados para su ubicación…"

# Clear the Times New Roman formatting from column A (rows 61-66), matching the rest of the sheet
$ws.Range("A61:A66").ClearFormats()
